# Auto-generated edit script: updates crypto price/volume table
# to reflect the refreshed data pulled on Fri May 12 22:46:13 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.824.69'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '1.811.78'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D5").Value = '309.18'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("E7").Value = '  +2.53%  '
$ws.Range("D8").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D8").Value = '0.3717'
$ws.Range("E8").Value = '  +3.30%  '
$ws.Range("D9").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D9").Value = '0.07250'
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D10").Value = '0.8682'
$ws.Range("E10").Value = '  +2.77%  '
$ws.Range("D11").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D11").Value = '20.86'
$ws.Range("E11").Value = '  +3.01%  '
$ws.Range("D12").Value = '1.986.66'
$ws.Range("E12").Value = '  +9.05%  '
$ws.Range("D13").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D13").Value = '6.670'
$ws.Range("E13").Value = '  +4.66%  '
$ws.Range("D14").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D14").Value = '5.363'
$ws.Range("E14").Value = '  +1.30%  '
$ws.Range("D15").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D15").Value = '0.06921'
$ws.Range("E15").Value = '  +2.20%  '
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D18").Value = '0.000008931'
$ws.Range("E18").Value = '  +2.58%  '
$ws.Range("D19").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D20").Value = '15.23'
$ws.Range("E20").Value = '  +1.49%  '
$ws.Range("D21").Value = '26.863.10'
$ws.Range("E21").Value = '  -0.94%  '
$ws.Range("D22").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D22").Value = '5.214'
$ws.Range("E22").Value = '  +2.83%  '
$ws.Range("D23").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D23").Value = '11.17'
$ws.Range("E23").Value = '  +1.20%  '
$ws.Range("D24").Value = '2.219.40'
$ws.Range("E24").Value = '  +8.50%  '
$ws.Range("D25").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D25").Value = '153.67'
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("D26").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D26").Value = '1.873'
$ws.Range("E26").Value = '  -2.74%  '
$ws.Range("D27").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D27").Value = '18.34'
$ws.Range("E27").Value = '  +1.10%  '
$ws.Range("D28").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D28").Value = '5.220'
$ws.Range("E28").Value = '  +3.86%  '
$ws.Range("D29").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D29").Value = '1.907'
$ws.Range("E29").Value = '  +15.25%  '
$ws.Range("D30").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D30").Value = '115.43'
$ws.Range("E30").Value = '  +1.80%  '
$ws.Range("D31").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D31").Value = '0.08944'
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("D32").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D32").Value = '0.7584'
$ws.Range("E32").Value = '  +3.53%  '
$ws.Range("D33").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D33").Value = '1.172'
$ws.Range("E33").Value = '  +7.10%  '
$ws.Range("D34").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D34").Value = '4.445'
$ws.Range("E34").Value = '  +2.34%  '
$ws.Range("D35").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D35").Value = '2.790'
$ws.Range("E35").Value = '  -2.53%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D36").Value = '1.138'
$ws.Range("E36").Value = '  +5.33%  '
$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D37").Value = '1.007'
$ws.Range("E37").Value = '  +0.53%  '
$ws.Range("D38").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D38").Value = '0.05233'
$ws.Range("E38").Value = '  +1.59%  '
$ws.Range("D39").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D39").Value = '0.01929'
$ws.Range("E39").Value = '  +1.41%  '
$ws.Range("D40").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D40").Value = '0.5092'
$ws.Range("E40").Value = '  +2.10%  '
$ws.Range("D41").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D41").Value = '0.1652'
$ws.Range("E41").Value = '  +1.25%  '
$ws.Range("D42").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D42").Value = '2.659'
$ws.Range("E42").Value = '  +1.58%  '
$ws.Range("D43").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D43").Value = '6.560'
$ws.Range("E43").Value = '  +9.88%  '
$ws.Range("D44").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D44").Value = '8.293'
$ws.Range("E44").Value = '  +2.60%  '
$ws.Range("D45").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D45").Value = '106.61'
$ws.Range("E45").Value = '  +1.42%  '
$ws.Range("D46").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D46").Value = '10.41'
$ws.Range("E46").Value = '  +2.04%  '
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("D48").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D48").Value = '1.659'
$ws.Range("E48").Value = '  +3.48%  '
$ws.Range("D49").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D49").Value = '0.4561'
$ws.Range("E49").Value = '  +0.32%  '
$ws.Range("D50").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D50").Value = '0.06272'
$ws.Range("E50").Value = '  -0.56%  '
$ws.Range("D51").NumberFormat = "@"  # keep as text, avoid numeric auto-conversion
$ws.Range("D51").Value = '1.811'
$ws.Range("E51").Value = '  +5.33%  '
